{"js": "// Replace the date line and the 25 two-digit multiplication problems\n// in the practice-sheet table with the new values from the commit.\nconst replacements = [\n  [\"2025-05-05 Monday\", \"2025-05-06 Tuesday\"],\n  [\"83\u00d794=\", \"34\u00d732=\"],\n  [\"20\u00d780=\", \"79\u00d725=\"],\n  [\"65\u00d780=\", \"37\u00d784=\"],\n  [\"63\u00d741=\", \"54\u00d761=\"],\n  [\"34\u00d798=\", \"65\u00d795=\"],\n  [\"84\u00d772=\", \"97\u00d757=\"],\n  [\"12\u00d752=\", \"63\u00d721=\"],\n  [\"45\u00d784=\", \"49\u00d727=\"],\n  [\"98\u00d790=\", \"85\u00d715=\"],\n  [\"29\u00d748=\", \"98\u00d787=\"],\n  [\"36\u00d760=\", \"87\u00d714=\"],\n  [\"97\u00d734=\", \"55\u00d711=\"],\n  [\"25\u00d782=\", \"72\u00d787=\"],\n  [\"65\u00d758=\", \"49\u00d725=\"],\n  [\"82\u00d763=\", \"46\u00d752=\"],\n  [\"78\u00d784=\", \"92\u00d713=\"],\n  [\"52\u00d792=\", \"46\u00d765=\"],\n  [\"39\u00d795=\", \"15\u00d763=\"],\n  [\"69\u00d760=\", \"58\u00d760=\"],\n  [\"22\u00d799=\", \"19\u00d785=\"],\n  [\"11\u00d725=\", \"54\u00d741=\"],\n  [\"63\u00d772=\", \"22\u00d779=\"],\n  [\"21\u00d770=\", \"12\u00d745=\"],\n  [\"49\u00d767=\", \"55\u00d782=\"],\n  [\"92\u00d755=\", \"20\u00d728=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit multiplication problems\n# in the practice-sheet table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-05 Monday\", \"2025-05-06 Tuesday\"),\n    @(\"83\u00d794=\", \"34\u00d732=\"),\n    @(\"20\u00d780=\", \"79\u00d725=\"),\n    @(\"65\u00d780=\", \"37\u00d784=\"),\n    @(\"63\u00d741=\", \"54\u00d761=\"),\n    @(\"34\u00d798=\", \"65\u00d795=\"),\n    @(\"84\u00d772=\", \"97\u00d757=\"),\n    @(\"12\u00d752=\", \"63\u00d721=\"),\n    @(\"45\u00d784=\", \"49\u00d727=\"),\n    @(\"98\u00d790=\", \"85\u00d715=\"),\n    @(\"29\u00d748=\", \"98\u00d787=\"),\n    @(\"36\u00d760=\", \"87\u00d714=\"),\n    @(\"97\u00d734=\", \"55\u00d711=\"),\n    @(\"25\u00d782=\", \"72\u00d787=\"),\n    @(\"65\u00d758=\", \"49\u00d725=\"),\n    @(\"82\u00d763=\", \"46\u00d752=\"),\n    @(\"78\u00d784=\", \"92\u00d713=\"),\n    @(\"52\u00d792=\", \"46\u00d765=\"),\n    @(\"39\u00d795=\", \"15\u00d763=\"),\n    @(\"69\u00d760=\", \"58\u00d760=\"),\n    @(\"22\u00d799=\", \"19\u00d785=\"),\n    @(\"11\u00d725=\", \"54\u00d741=\"),\n    @(\"63\u00d772=\", \"22\u00d779=\"),\n    @(\"21\u00d770=\", \"12\u00d745=\"),\n    @(\"49\u00d767=\", \"55\u00d782=\"),\n    @(\"92\u00d755=\", \"20\u00d728=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
